$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C:G (est, l95, u95, p, interaction), rows 2-21
$values = @(
    @(4.54, 1.81, 7.28, 0.00124, 0.0658),
    @(3.74, 0.54, 6.94, 0.02236, 0.02741),
    @(4.36, 1.33, 7.39, 0.00497, 0.00396),
    @(0.64, -3.44, 4.72, 0.75688, 0.50157),
    @(5.3, 2.64, 7.96, 0.00012, 0.03499),
    @(3.86, 1.05, 6.67, 0.00736, 0.38106),
    @(3.71, 0.93, 6.5, 0.00931, 0.36459),
    @(-0.19, -3.43, 3.05, 0.90706, 0.19079),
    @(3.11, 0.23, 6, 0.03436, 0.50781),
    @(-1.38, -5.31, 2.56, 0.4917, 0.52876),
    @(1.72, -0.94, 4.39, 0.2033, 0.0658),
    @(-0.74, -4.13, 2.64, 0.6648, 0.02741),
    @(-0.22, -2.9, 2.45, 0.86826, 0.00396),
    @(0.24, -3.76, 4.24, 0.90569, 0.50157),
    @(2.5, -0.19, 5.18, 0.06812, 0.03499),
    @(2.69, -0.33, 5.72, 0.08071, 0.38106),
    @(5.18, 2.09, 8.27, 0.00116, 0.36459),
    @(2.89, 0.15, 5.62, 0.03851, 0.19079),
    @(1.55, -1.96, 5.06, 0.38464, 0.50781),
    @(-2.04, -5.63, 1.54, 0.26206, 0.52876)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $values[$i]
    $ws.Cells.Item($row, 3).Value = $rowValues[0]
    $ws.Cells.Item($row, 4).Value = $rowValues[1]
    $ws.Cells.Item($row, 5).Value = $rowValues[2]
    $ws.Cells.Item($row, 6).Value = $rowValues[3]
    $ws.Cells.Item($row, 7).Value = $rowValues[4]
}
